# Poland II Liga workbook update (16-06-2024)
# The underlying data rows in several blocks were out of order relative to
# their match id (column B). This script restores the correct row contents
# (columns B..AD) for those rows while leaving column A (the running index)
# untouched, by reading the current contents of every row in a block and
# then rewriting each row with the contents that belong there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Reorder-Rows {
    param($Worksheet, $Mapping)

    # Snapshot the current B:AD contents of every row involved.
    $snapshot = @{}
    foreach ($row in $Mapping.Keys) {
        $rng = $Worksheet.Range("B$row" + ":AD$row")
        $snapshot[$row] = $rng.Value()
    }

    # Write back the row that belongs at each position.
    foreach ($row in $Mapping.Keys) {
        $src = $Mapping[$row]
        $rng = $Worksheet.Range("B$row" + ":AD$row")
        $rng.Value = $snapshot[$src]
    }
}

# Block 1: rows 3-8 (ids 5226043, 5228058, 5185127, 5224889, 5229462, 5226044)
$map1 = @{
    3 = 4
    4 = 3
    5 = 8
    6 = 7
    7 = 5
    8 = 6
}
Reorder-Rows $ws $map1

# Block 2: rows 65-66
$map2 = @{
    65 = 66
    66 = 65
}
Reorder-Rows $ws $map2

# Block 3: rows 233-234
$map3 = @{
    233 = 234
    234 = 233
}
Reorder-Rows $ws $map3

# Block 4: rows 271-273
$map4 = @{
    271 = 273
    272 = 271
    273 = 272
}
Reorder-Rows $ws $map4

# Block 5: rows 306-312 (307 keeps its own content)
$map5 = @{
    306 = 312
    307 = 307
    308 = 306
    309 = 311
    310 = 308
    311 = 309
    312 = 310
}
Reorder-Rows $ws $map5
